$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  0.30423342464539616,
  -0.030381441205175719,
  -0.0039999999953463927,
  -0.0079999999918545228,
  -0.0029999999950103273,
  -0.0019999999944477764,
  -0.009999999988560937,
  0.0044535460831327534,
  -0.0019999999941759938,
  -0.0019999999940285562,
  -0.0029999999932952548,
  -0.003499999992883307,
  -0.0034999999928011505,
  -0.00799999998950085,
  0.013857821434005047,
  -0.0019999999939219748,
  -0.0019999999938358215,
  -0.0039999999923390206,
  -0.0039999999967443856,
  -0.007596566054196785,
  -0.0039999999964717148,
  -0.0039999999964370758,
  -0.069024928240546402,
  -0.019999999982576178,
  -0.019999999982356798,
  -0.0024999999939012696,
  -0.0024999999935553241,
  -0.0019999999926731959,
  -0.0069999999880003827,
  -0.059999999948468385,
  -0.0069999999871566132,
  -0.0099999999848847665,
  -0.0039999999892899041
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
